# ------------------------------------------------------------------
# Update "上海-漫展信息.xlsx" per commit: "Update gh-pages to output
# generated at 456a3b4" -- refreshes visitor-interest counters (F
# column) across sheets, flips one local-life listing to "已停售",
# and regenerates the combined "全部类型" sheet (one stale row
# removed, everything below shifted up, counters refreshed).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibitions) -- refresh "想去人数" (F column) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 590
$ws1.Range("F4").Value  = 6369
$ws1.Range("F7").Value  = 71
$ws1.Range("F8").Value  = 318
$ws1.Range("F9").Value  = 194
$ws1.Range("F10").Value = 11
$ws1.Range("F11").Value = 693
$ws1.Range("F12").Value = 1172
$ws1.Range("F14").Value = 419
$ws1.Range("F17").Value = 1417
$ws1.Range("F19").Value = 380
$ws1.Range("F20").Value = 397
$ws1.Range("F22").Value = 1069
$ws1.Range("F23").Value = 138
$ws1.Range("F24").Value = 2208
$ws1.Range("F25").Value = 254
$ws1.Range("F26").Value = 95
$ws1.Range("F27").Value = 394
$ws1.Range("F28").Value = 58
$ws1.Range("F29").Value = 3565
$ws1.Range("F30").Value = 45
$ws1.Range("F31").Value = 632

# ---- Sheet 2: 演出 (Performances) -- refresh F column ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value  = 358
$ws2.Range("F4").Value  = 175
$ws2.Range("F8").Value  = 706
$ws2.Range("F11").Value = 1017
$ws2.Range("F13").Value = 101
$ws2.Range("F18").Value = 377
$ws2.Range("F24").Value = 191
$ws2.Range("F32").Value = 1657
$ws2.Range("F33").Value = 23

# ---- Sheet 3: 本地生活 (Local life) -- refresh F column + sold-out flag ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F9").Value  = 129
$ws3.Range("F11").Value = 768
$ws3.Range("G2").Value  = "已停售"

# ---- Sheet 4: 全部类型 (All types, combined/regenerated listing) ----
# The first listing (2023.01.12 上海·日漫咖啡体验) drops out of the
# combined sheet entirely; every row below shifts up by one, and the
# dimension shrinks from I50 to I49.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows("2:2").Delete()

# Re-number the running index column (A) which holds literal numbers,
# not a formula, so it must be fixed up after the shift.
for ($r = 2; $r -le 49; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

# Apply the same counter refreshes as above, at their new (shifted) row numbers.
$ws4.Range("F6").Value  = 129
$ws4.Range("F7").Value  = 768
$ws4.Range("F8").Value  = 590
$ws4.Range("F10").Value = 6369
$ws4.Range("F14").Value = 706
$ws4.Range("F15").Value = 71
$ws4.Range("F16").Value = 318
$ws4.Range("F17").Value = 194
$ws4.Range("F18").Value = 693
$ws4.Range("F19").Value = 101
$ws4.Range("F21").Value = 1172
$ws4.Range("F23").Value = 419
$ws4.Range("F25").Value = 377
$ws4.Range("F27").Value = 1417
$ws4.Range("F30").Value = 380
$ws4.Range("F31").Value = 397
$ws4.Range("F33").Value = 191
$ws4.Range("F35").Value = 1069
$ws4.Range("F36").Value = 138
$ws4.Range("F38").Value = 2208
$ws4.Range("F40").Value = 1657
$ws4.Range("F41").Value = 254
$ws4.Range("F42").Value = 95
$ws4.Range("F43").Value = 394
$ws4.Range("F44").Value = 58
$ws4.Range("F45").Value = 3565
$ws4.Range("F48").Value = 45
$ws4.Range("F49").Value = 632
